$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP1912")

# Header updates
$ws.Range("A2").Value = "Última actualización: 13:33:42"
$ws.Range("A3").Value = "Total filas: 258"

# Cell-level corrections for re-ordered/updated rows
# Row 53
$ws.Range("C53").Value = "16_SANTA ANA"

# Row 54
$ws.Range("C54").Value = "11_ETCHEVERRY"

# Row 70
$ws.Range("A70").Value = "06:55:02"
$ws.Range("C70").Value = "23_HERNANDEZ"
$ws.Range("D70").Value = 68

# Row 71
$ws.Range("A71").Value = "07:19:29"
$ws.Range("C71").Value = "11_ETCHEVERRY"
$ws.Range("D71").Value = 44

# Row 103
$ws.Range("A103").Value = "07:50:16"
$ws.Range("C103").Value = "17_ROMERO"
$ws.Range("D103").Value = 93

# Row 104
$ws.Range("A104").Value = "08:32:09"
$ws.Range("C104").Value = "16_SANTA ANA"
$ws.Range("D104").Value = 51

# Row 105
$ws.Range("A105").Value = "08:02:22"
$ws.Range("C105").Value = "11_ETCHEVERRY"
$ws.Range("D105").Value = 81

# Row 134
$ws.Range("A134").Value = "09:35:26"
$ws.Range("C134").Value = "23_HERNANDEZ"
$ws.Range("D134").Value = 59

# Row 135
$ws.Range("A135").Value = "10:29:57"
$ws.Range("C135").Value = "16_SANTA ANA"
$ws.Range("D135").Value = 5

# Row 150
$ws.Range("A150").Value = "10:59:49"
$ws.Range("C150").Value = "23_HERNANDEZ"
$ws.Range("D150").Value = 7

# Row 151
$ws.Range("A151").Value = "09:35:26"
$ws.Range("C151").Value = "16_P MOR-167 Y 521"
$ws.Range("D151").Value = 91

# Row 215
$ws.Range("A215").Value = "12:21:08"
$ws.Range("C215").Value = "16_P MOR-SANTA ANA"
$ws.Range("D215").Value = 46

# Row 216
$ws.Range("A216").Value = "11:30:45"
$ws.Range("C216").Value = "10_OLMOS"
$ws.Range("D216").Value = 97

# Row 231
$ws.Range("A231").Value = "13:33:42"
$ws.Range("C231").Value = "16_SANTA ANA"
$ws.Range("D231").Value = 0

# Row 232
$ws.Range("A232").Value = "12:59:47"
$ws.Range("B232").Value = "13:33"
$ws.Range("C232").Value = "10_OLMOS"
$ws.Range("D232").Value = 34

# Row 233
$ws.Range("A233").Value = "13:33:42"
$ws.Range("B233").Value = "13:34"
$ws.Range("C233").Value = "16_SANTA ANA"
$ws.Range("D233").Value = 1

# Row 234
$ws.Range("A234").Value = "13:33:42"
$ws.Range("B234").Value = "13:34"
$ws.Range("C234").Value = "23_HERNANDEZ"
$ws.Range("D234").Value = 1

# Row 235
$ws.Range("A235").Value = "11:56:55"
$ws.Range("B235").Value = "13:36"
$ws.Range("C235").Value = "15_ABASTO"
$ws.Range("D235").Value = 100

# Row 236
$ws.Range("A236").Value = "13:33:42"
$ws.Range("B236").Value = "13:38"
$ws.Range("C236").Value = "14_ABASTO"
$ws.Range("D236").Value = 5

# Row 237
$ws.Range("B237").Value = "13:46"
$ws.Range("C237").Value = "17_ROMERO"
$ws.Range("D237").Value = 110

# Row 238
$ws.Range("A238").Value = "12:59:47"
$ws.Range("B238").Value = "13:50"
$ws.Range("C238").Value = "11_ETCHEVERRY"
$ws.Range("D238").Value = 51

# Row 239
$ws.Range("A239").Value = "11:56:55"
$ws.Range("B239").Value = "13:50"
$ws.Range("C239").Value = "215A_EL PATO"
$ws.Range("D239").Value = 114

# Row 240
$ws.Range("A240").Value = "12:21:08"
$ws.Range("B240").Value = "13:51"
$ws.Range("C240").Value = "215A_EL PATO"
$ws.Range("D240").Value = 90

# Row 241
$ws.Range("A241").Value = "11:56:55"
$ws.Range("B241").Value = "13:55"
$ws.Range("C241").Value = "225_GOMEZ"
$ws.Range("D241").Value = 119

# Row 242
$ws.Range("B242").Value = "13:56"
$ws.Range("C242").Value = "225_GOMEZ"
$ws.Range("D242").Value = 95

# Row 243
$ws.Range("A243").Value = "12:59:47"
$ws.Range("B243").Value = "13:56"
$ws.Range("C243").Value = "16_P MOR-167 Y 521"
$ws.Range("D243").Value = 57

# Row 244
$ws.Range("A244").Value = "12:47:27"
$ws.Range("B244").Value = "13:58"
$ws.Range("C244").Value = "16_P MOR-167 Y 521"
$ws.Range("D244").Value = 71

# Row 245
$ws.Range("A245").Value = "12:21:08"
$ws.Range("B245").Value = "14:00"
$ws.Range("C245").Value = "16_P MOR-167 Y 521"
$ws.Range("D245").Value = 99

# Row 246
$ws.Range("B246").Value = "14:04"
$ws.Range("C246").Value = "17_ROMERO"
$ws.Range("D246").Value = 103

# Row 247
$ws.Range("A247").Value = "13:33:42"
$ws.Range("B247").Value = "14:04"
$ws.Range("C247").Value = "23_HERNANDEZ"
$ws.Range("D247").Value = 31

# Row 248
$ws.Range("B248").Value = "14:08"
$ws.Range("C248").Value = "23_HERNANDEZ"
$ws.Range("D248").Value = 107

# Row 249
$ws.Range("A249").Value = "12:59:47"
$ws.Range("B249").Value = "14:11"
$ws.Range("C249").Value = "23_HERNANDEZ"
$ws.Range("D249").Value = 72

# Row 250
$ws.Range("A250").Value = "13:33:42"
$ws.Range("B250").Value = "14:12"
$ws.Range("C250").Value = "15_ABASTO"
$ws.Range("D250").Value = 39

# Row 251
$ws.Range("A251").Value = "12:47:27"
$ws.Range("B251").Value = "14:16"
$ws.Range("C251").Value = "27_EL RETIRO"
$ws.Range("D251").Value = 89

# New rows appended at the end (252-263)
# Row 252
$ws.Range("A252").Value = "12:21:08"
$ws.Range("B252").Value = "14:17"
$ws.Range("C252").Value = "27_EL RETIRO"
$ws.Range("D252").Value = 116
$ws.Range("E252").Value = "LP1912"

# Row 253
$ws.Range("A253").Value = "12:59:47"
$ws.Range("B253").Value = "14:19"
$ws.Range("C253").Value = "215C_EL PATO"
$ws.Range("D253").Value = 80
$ws.Range("E253").Value = "LP1912"

# Row 254
$ws.Range("A254").Value = "12:21:08"
$ws.Range("B254").Value = "14:20"
$ws.Range("C254").Value = "215C_EL PATO"
$ws.Range("D254").Value = 119
$ws.Range("E254").Value = "LP1912"

# Row 255
$ws.Range("A255").Value = "12:47:27"
$ws.Range("B255").Value = "14:21"
$ws.Range("C255").Value = "26_HERNANDEZ"
$ws.Range("D255").Value = 94
$ws.Range("E255").Value = "LP1912"

# Row 256
$ws.Range("A256").Value = "13:33:42"
$ws.Range("B256").Value = "14:44"
$ws.Range("C256").Value = "14_ABASTO"
$ws.Range("D256").Value = 71
$ws.Range("E256").Value = "LP1912"

# Row 257
$ws.Range("A257").Value = "12:47:27"
$ws.Range("B257").Value = "14:45"
$ws.Range("C257").Value = "14_ABASTO"
$ws.Range("D257").Value = 118
$ws.Range("E257").Value = "LP1912"

# Row 258
$ws.Range("A258").Value = "12:59:47"
$ws.Range("B258").Value = "14:56"
$ws.Range("C258").Value = "16_P MOR-SANTA ANA"
$ws.Range("D258").Value = 117
$ws.Range("E258").Value = "LP1912"

# Row 259
$ws.Range("A259").Value = "12:59:47"
$ws.Range("B259").Value = "14:58"
$ws.Range("C259").Value = "215B_EL PATO"
$ws.Range("D259").Value = 119
$ws.Range("E259").Value = "LP1912"

# Row 260
$ws.Range("A260").Value = "13:33:42"
$ws.Range("B260").Value = "15:00"
$ws.Range("C260").Value = "81_EL PELIGRO"
$ws.Range("D260").Value = 87
$ws.Range("E260").Value = "LP1912"

# Row 261
$ws.Range("A261").Value = "13:33:42"
$ws.Range("B261").Value = "15:05"
$ws.Range("C261").Value = "10_OLMOS"
$ws.Range("D261").Value = 92
$ws.Range("E261").Value = "LP1912"

# Row 262
$ws.Range("A262").Value = "13:33:42"
$ws.Range("B262").Value = "15:13"
$ws.Range("C262").Value = "11_ETCHEVERRY"
$ws.Range("D262").Value = 100
$ws.Range("E262").Value = "LP1912"

# Row 263
$ws.Range("A263").Value = "13:33:42"
$ws.Range("B263").Value = "15:17"
$ws.Range("C263").Value = "26_HERNANDEZ"
$ws.Range("D263").Value = 104
$ws.Range("E263").Value = "LP1912"

# Sibling sheets: refresh timestamp only
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 13:33:42"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 13:33:42"